# Adds a new row of PIM-user test data (page object model) to the
# AddPIMUserData sheet, and normalises the stale row/column selection
# markers left over on all three sheets.

$wb = $excel.ActiveWorkbook

# --- LoginData: selection sqref fix (was "1:1", should be "B1") ---
$wsLogin = $wb.Worksheets.Item("LoginData")
$wsLogin.Range("B1").Select()

# --- AddUserData: selection sqref fix (was "1:1", should be "A1") ---
$wsAddUser = $wb.Worksheets.Item("AddUserData")
$wsAddUser.Range("A1").Select()

# --- AddPIMUserData: new page-object-model row + refreshed header row ---
$ws = $wb.Worksheets.Item("AddPIMUserData")

# Drop the two existing mailto hyperlinks (on E1 / F1) before rebuilding
# the row -- they get re-created further right, at H1/I1.
$ws.Range("E1").Hyperlinks.Delete()

# Row 1: existing Admin/admin123 login stay in A1/B1; C1:G1 become the new
# "Michale Joseph Martin" user fields, H1/I1 become the two hyperlink cells.
$ws.Range("C1").Value = "Michale"
$ws.Range("D1").Value = "Joseph"
$ws.Range("E1").Value = "Martin"
$ws.Range("F1").Value = 20003
$ws.Range("G1").Value = "mjmartin"
$ws.Hyperlinks.Add($ws.Range("H1"), "mailto:Test@123", [Type]::Missing, [Type]::Missing, "Test@123")
$ws.Range("H1").Style = $ws.Range("A1").Style
$ws.Hyperlinks.Add($ws.Range("I1"), "mailto:Test@123", [Type]::Missing, [Type]::Missing, "Test@123")
$ws.Range("I1").Style = $ws.Range("A1").Style

# Row 2: brand-new "Samuel Jack Todd" user row.
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("C2").Value = "Samuel"
$ws.Range("D2").Value = "Jack"
$ws.Range("E2").Value = "Todd"
$ws.Range("F2").Value = 2002
$ws.Range("G2").Value = "Sjtodd"
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:Test@123", [Type]::Missing, [Type]::Missing, "Test@123")
$ws.Range("H2").Style = $ws.Range("A1").Style
$ws.Hyperlinks.Add($ws.Range("I2"), "mailto:Test@123", [Type]::Missing, [Type]::Missing, "Test@123")
$ws.Range("I2").Style = $ws.Range("A1").Style

# Active selection ends on G2 (the last new data cell) on the tab-selected sheet.
$ws.Range("G2").Select()
